$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: year, S0, S1, S2
$newRows = @(
    @(2035, 0.8099999999999999, 0.15, 0.04),
    @(2040, 0.8099999999999999, 0.15, 0.04),
    @(2045, 0.8099999999999999, 0.15, 0.04),
    @(2050, 0.8099999999999999, 0.15, 0.04)
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
